$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 1 de Septiembre de 2020 a las 07:33"

# Pakistan (row 19) data refresh
$ws.Range("B19").Value = 296149
$ws.Range("C19").Value = 300
$ws.Range("D19").Value = 280970
$ws.Range("E19").Value = 8881
$ws.Range("G19").Value = 4
$ws.Range("H19").Value = 6298

# Filipinas (row 25) data refresh
$ws.Range("D25").Value = 157559
$ws.Range("E25").Value = 59697
$ws.Range("H25").Value = 3563

# Israel moves above Bolivia in the country ranking (rows 30-31 swap with
# updated figures for Israel, Bolivia keeps its previous totals)
$ws.Range("A30").Value = "Israel"
$ws.Range("B30").Value = 117030
$ws.Range("C30").Value = 434
$ws.Range("D30").Value = 95589
$ws.Range("E30").Value = 20502
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 939

$ws.Range("A31").Value = "Bolivia"
$ws.Range("B31").Value = 116598
$ws.Range("C31").Value = 630
$ws.Range("D31").Value = 60408
$ws.Range("E31").Value = 51163
$ws.Range("G31").Value = 61
$ws.Range("H31").Value = 5027

# Kirguistan (row 59) data refresh
$ws.Range("B59").Value = 43958
$ws.Range("C59").Value = 60
$ws.Range("D59").Value = 38649
$ws.Range("E59").Value = 4250
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 1059

# Tailandia (row 124) data refresh
$ws.Range("B124").Value = 3417
$ws.Range("C124").Value = 5
$ws.Range("D124").Value = 3274
$ws.Range("E124").Value = 85

# Curazao (row 196) data refresh
$ws.Range("B196").Value = 69
$ws.Range("C196").Value = 1
$ws.Range("E196").Value = 33
